# Weekly update: insert a new price record for Berenjena (Vega Monumental
# Concepción) before the current row 70, shifting all subsequent rows
# (old 70-114) down by one to (71-115).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 70; Excel shifts rows 70..114 down to 71..115
# and carries column formatting (e.g. the date style on column D) along.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new record's data.
$ws.Cells.Item(70, 1).Value2  = 11
$ws.Cells.Item(70, 2).Value   = "Vega Monumental Concepción"
$ws.Cells.Item(70, 3).Value   = "Bíobío"
$ws.Cells.Item(70, 4).Value2  = 44806
$ws.Cells.Item(70, 5).Value2  = 8
$ws.Cells.Item(70, 6).Value2  = 100112001
$ws.Cells.Item(70, 7).Value   = "Berenjena"
$ws.Cells.Item(70, 8).Value   = "Sin especificar"
$ws.Cells.Item(70, 9).Value   = "Primera"
$ws.Cells.Item(70, 10).Value2 = 220
$ws.Cells.Item(70, 11).Value2 = 12000
$ws.Cells.Item(70, 12).Value2 = 13000
$ws.Cells.Item(70, 13).Value2 = 12545
$ws.Cells.Item(70, 14).Value  = "$/caja 60 unidades"
$ws.Cells.Item(70, 15).Value  = "Región de Arica y Parinacota"
$ws.Cells.Item(70, 16).Value2 = 209
$ws.Cells.Item(70, 17).Value2 = 60
$ws.Cells.Item(70, 18).Value  = "Hortaliza"
